$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (all source cells are plain text / inlineStr)
$updates = @{
    'D2' = '243.51'
    'E2' = '-0.56%'
    'G2' = '4'
    'D3' = '30.24'
    'E3' = '13.89%'
    'G3' = '4'
    'D4' = '5.140'
    'E4' = '0.32%'
    'G4' = '4'
    'D5' = '0.05670'
    'E5' = '1.45%'
    'G5' = '4'
    'D6' = '6.523'
    'E6' = '0.78%'
    'G6' = '4'
    'D7' = '0.8319'
    'E7' = '1.83%'
    'G7' = '4'
    'D8' = '0.8614'
    'E8' = '3.19%'
    'G8' = '4'
    'D9' = '0.01009'
    'E9' = '1,581.93%'
    'G9' = '4'
    'D10' = '0.1332'
    'E10' = '-0.27%'
    'G10' = '4'
    'D11' = '0.06917'
    'G11' = '4'
    'D12' = '0.02869'
    'E12' = '-0.58%'
    'G12' = '4'
    'D13' = '0.09373'
    'E13' = '-0.16%'
    'G13' = '4'
    'E14' = '-0.42%'
    'G14' = '4'
    'D15' = '0.04154'
    'E15' = '-9.63%'
    'G15' = '4'
    'D16' = '0.006001'
    'E16' = '-3.32%'
    'G16' = '4'
    'D17' = '3.520'
    'E17' = '-3.53%'
    'G17' = '4'
    'D18' = '3.022'
    'E18' = '-0.52%'
    'G18' = '4'
    'D19' = '2.131'
    'E19' = '-2.38%'
    'G19' = '4'
    'D20' = '0.3151'
    'E20' = '1.27%'
    'G20' = '4'
    'D21' = '0.03274'
    'E21' = '4.69%'
    'G21' = '4'
    'E22' = '-0.31%'
    'G22' = '4'
    'D23' = '3.625'
    'E23' = '-3.64%'
    'G23' = '4'
    'E24' = '-0.12%'
    'G24' = '4'
    'D25' = '0.001211'
    'E25' = '-2.74%'
    'G25' = '4'
    'D26' = '0.004450'
    'E26' = '-1.17%'
    'G26' = '4'
    'D27' = '0.0001179'
    'E27' = '22.76%'
    'G27' = '4'
    'D28' = '0.0001396'
    'E28' = '0.21%'
    'G28' = '4'
    'G29' = '4'
    'G30' = '4'
    'G31' = '4'
    'G32' = '4'
    'G33' = '4'
    'G34' = '4'
    'G35' = '4'
    'G36' = '4'
    'G37' = '4'
    'G38' = '4'
    'G39' = '4'
    'D40' = '0.03713'
    'E40' = '1.93%'
    'G40' = '4'
    'D41' = '0.005824'
    'E41' = '-5.94%'
    'G41' = '4'
    'D42' = '0.1055'
    'E42' = '0.46%'
    'G42' = '4'
    'D43' = '0.002309'
    'E43' = '-3.82%'
    'G43' = '4'
    'D44' = '0.009746'
    'E44' = '10.18%'
    'G44' = '4'
    'D45' = '0.00005093'
    'E45' = '-4.58%'
    'G45' = '4'
    'D46' = '0.00000000749'
    'E46' = '-0.12%'
    'G46' = '4'
    'D47' = '0.09990'
    'E47' = '-30.63%'
    'G47' = '4'
    'D48' = '0.002878'
    'E48' = '23.14%'
    'G48' = '4'
    'D49' = '0.00002098'
    'E49' = '-0.12%'
    'G49' = '4'
    'D50' = '0.0001998'
    'E50' = '-0.12%'
    'G50' = '4'
    'G51' = '4'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text entry so numeric-/percent-looking strings are not
    # auto-converted to numbers, then restore the default cell style
    # so no stray formatting is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
